$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2587845
$ws.Range("J17").Value = 2587845
$ws.Range("L17").Value = 7763535
$ws.Range("N17").Value = -7763871
$ws.Range("H20").Value = 4900
$ws.Range("I20").Value = 4900
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 4900
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -4670
$ws.Range("H35").Value = 4900
$ws.Range("I35").Value = 4900
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4900
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -4521
$ws.Range("H64").Value = 6918.1816
$ws.Range("I64").Value = 8012.5
$ws.Range("K64").Value = 8012.5
$ws.Range("M64").Value = -7764.5
$ws.Range("H67").Value = 6918.1816
$ws.Range("I67").Value = 8012.5
$ws.Range("K67").Value = 8012.5
$ws.Range("M67").Value = -7154.5
$ws.Range("H69").Value = 5262.1875
$ws.Range("I69").Value = 4640
$ws.Range("J69").Value = 5545
$ws.Range("K69").Value = 13920
$ws.Range("L69").Value = 16635
$ws.Range("M69").Value = -13046
$ws.Range("N69").Value = -18383
$ws.Range("H72").Value = 5262.1875
$ws.Range("I72").Value = 4640
$ws.Range("J72").Value = 5545
$ws.Range("K72").Value = 41760
$ws.Range("L72").Value = 49905
$ws.Range("M72").Value = -37392
$ws.Range("N72").Value = -58641
$ws.Range("H74").Value = 4916.3335
$ws.Range("I74").Value = 3740
$ws.Range("J74").Value = 5756.5713
$ws.Range("K74").Value = 3740
$ws.Range("L74").Value = 5756.5713
$ws.Range("M74").Value = -2804
$ws.Range("N74").Value = -7628.5713
$ws.Range("H77").Value = 4916.3335
$ws.Range("I77").Value = 3740
$ws.Range("J77").Value = 5756.5713
$ws.Range("K77").Value = 18700
$ws.Range("L77").Value = 28782.8565
$ws.Range("M77").Value = -14020
$ws.Range("N77").Value = -38142.85649999999
$ws.Range("H80").Value = 983.3333
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 975
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 2925
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4921
$ws.Range("H83").Value = 983.3333
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 975
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 8775
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -18759
$ws.Range("H88").Value = 9569.923000000001
$ws.Range("I88").Value = 1700
$ws.Range("J88").Value = 11000.818
$ws.Range("K88").Value = 1700
$ws.Range("L88").Value = 11000.818
$ws.Range("M88").Value = -1294
$ws.Range("N88").Value = -11812.818
$ws.Range("H91").Value = 9569.923000000001
$ws.Range("I91").Value = 1700
$ws.Range("J91").Value = 11000.818
$ws.Range("K91").Value = 1700
$ws.Range("L91").Value = 11000.818
$ws.Range("M91").Value = -296
$ws.Range("N91").Value = -13808.818
$ws.Range("H106").Value = 57973600
$ws.Range("I106").Value = 22224852
$ws.Range("K106").Value = 22224852
$ws.Range("M106").Value = -22224221
$ws.Range("H113").Value = 9617962
$ws.Range("I113").Value = 2792.5
$ws.Range("J113").Value = 125000000
$ws.Range("K113").Value = 2792.5
$ws.Range("L113").Value = 125000000
$ws.Range("M113").Value = 461.5
$ws.Range("N113").Value = -125006508
$ws.Range("H137").Value = 1483
$ws.Range("I137").Value = 1192.0333
$ws.Range("K137").Value = 3576.0999
$ws.Range("M137").Value = -1026.0999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2040.34
$ws.Range("I32").Value = 1753.1538
$ws.Range("J32").Value = 4944.1113
$ws.Range("K32").Value = 1753.1538
$ws.Range("L32").Value = 4944.1113
$ws.Range("M32").Value = -1466.1538
$ws.Range("N32").Value = -5518.1113
$ws.Range("H63").Value = 250007500
$ws.Range("I63").Value = 250007500
$ws.Range("K63").Value = 250007500
$ws.Range("M63").Value = -250006814
$ws.Range("H66").Value = 250007500
$ws.Range("I66").Value = 250007500
$ws.Range("K66").Value = 1250037500
$ws.Range("M66").Value = -1250034068
$ws.Range("H88").Value = 2399
$ws.Range("J88").Value = 2399
$ws.Range("L88").Value = 2399
$ws.Range("N88").Value = -3211
$ws.Range("H91").Value = 2399
$ws.Range("J91").Value = 2399
$ws.Range("L91").Value = 2399
$ws.Range("N91").Value = -5207
$ws.Range("H132").Value = 2531.3518
$ws.Range("I132").Value = 1803.8684
$ws.Range("J132").Value = 4259.125
$ws.Range("K132").Value = 5411.6052
$ws.Range("L132").Value = 12777.375
$ws.Range("M132").Value = -2881.6052
$ws.Range("N132").Value = -17837.375

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17545374
$ws.Range("I86").Value = 23810974
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 23810974
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -23809851
$ws.Range("N86").Value = -3946
$ws.Range("H89").Value = 17545374
$ws.Range("I89").Value = 23810974
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 119054870
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -119049254
$ws.Range("N89").Value = -19732
$ws.Range("H105").Value = 62501800
$ws.Range("I105").Value = 83334984
$ws.Range("K105").Value = 83334984
$ws.Range("M105").Value = -83333237

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6026605.5
$ws.Range("I31").Value = 1632.8478
$ws.Range("K31").Value = 1632.8478
$ws.Range("M31").Value = -1337.8478
$ws.Range("H34").Value = 6026605.5
$ws.Range("I34").Value = 1632.8478
$ws.Range("K34").Value = 1632.8478
$ws.Range("M34").Value = -1430.8478
$ws.Range("H58").Value = 3705015
$ws.Range("I58").Value = 5556284.5
$ws.Range("K58").Value = 5556284.5
$ws.Range("M58").Value = -5556081.5
$ws.Range("H74").Value = 19790.223
$ws.Range("J74").Value = 22103.375
$ws.Range("L74").Value = 22103.375
$ws.Range("N74").Value = -23851.375
$ws.Range("H77").Value = 19790.223
$ws.Range("J77").Value = 22103.375
$ws.Range("L77").Value = 66310.125
$ws.Range("N77").Value = -75046.125
$ws.Range("H132").Value = 3126525.5
$ws.Range("I132").Value = 3704797.2
$ws.Range("J132").Value = 3858.9
$ws.Range("K132").Value = 11114391.6
$ws.Range("L132").Value = 11576.7
$ws.Range("M132").Value = -11111861.6
$ws.Range("N132").Value = -16636.7
$ws.Range("H136").Value = 3705015
$ws.Range("I136").Value = 5556284.5
$ws.Range("K136").Value = 16668853.5
$ws.Range("M136").Value = -16666303.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2222790.8
$ws.Range("I113").Value = 2500586.5
$ws.Range("J113").Value = 1429088.4
$ws.Range("K113").Value = 7501759.5
$ws.Range("L113").Value = 4287265.199999999
$ws.Range("M113").Value = -7499589.5
$ws.Range("N113").Value = -4291605.199999999
$ws.Range("H132").Value = 22223622
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 27779278
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 250013502
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -250018562

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2706.25
$ws.Range("I80").Value = 2390.9092
$ws.Range("J80").Value = 3400
$ws.Range("K80").Value = 2390.9092
$ws.Range("L80").Value = 3400
$ws.Range("M80").Value = -1392.9092
$ws.Range("N80").Value = -5396
$ws.Range("H83").Value = 2706.25
$ws.Range("I83").Value = 2390.9092
$ws.Range("J83").Value = 3400
$ws.Range("K83").Value = 11954.546
$ws.Range("L83").Value = 17000
$ws.Range("M83").Value = -6962.546
$ws.Range("N83").Value = -26984
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 381647.38
$ws.Range("I82").Value = 501363.7
$ws.Range("K82").Value = 501363.7
$ws.Range("M82").Value = -501002.7
$ws.Range("H85").Value = 381647.38
$ws.Range("I85").Value = 501363.7
$ws.Range("K85").Value = 501363.7
$ws.Range("M85").Value = -500115.7
$ws.Range("H94").Value = 26466
$ws.Range("J94").Value = 26466
$ws.Range("L94").Value = 26466
$ws.Range("N94").Value = -27818
$ws.Range("H130").Value = 38500
$ws.Range("J130").Value = 38500
$ws.Range("L130").Value = 38500
$ws.Range("N130").Value = -48540
$ws.Range("H132").Value = 11183285
$ws.Range("I132").Value = 13893521
$ws.Range("J132").Value = 3562.25
$ws.Range("K132").Value = 41680563
$ws.Range("L132").Value = 10686.75
$ws.Range("M132").Value = -41678033
$ws.Range("N132").Value = -15746.75
$ws.Range("H136").Value = 5427.207
$ws.Range("I136").Value = 3716.087
$ws.Range("K136").Value = 11148.261
$ws.Range("M136").Value = -8598.261

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2400
$ws.Range("I62").Value = 2400
$ws.Range("K62").Value = 2400
$ws.Range("M62").Value = -1776
$ws.Range("H65").Value = 2400
$ws.Range("I65").Value = 2400
$ws.Range("K65").Value = 12000
$ws.Range("M65").Value = -8880
$ws.Range("H125").Value = 35357.5
$ws.Range("J125").Value = 35357.5
$ws.Range("L125").Value = 35357.5
$ws.Range("N125").Value = -45197.5
